$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh prepends a new pair of rows (Primera / Segunda) for the
# latest market date, pushing every existing data row down by two rows.
$ws.Rows.Item(2).EntireRow.Insert()
$ws.Rows.Item(2).EntireRow.Insert()

# The freshly inserted rows inherit the header row's bold/centered style;
# reset them back to the plain style used by the rest of the data rows.
$ws.Range("A2:R3").ClearFormats()

# New row 2: Primera quality for the new date (2022-05-18)
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44699
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("N2").Value = "$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 861
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"

# New row 3: Segunda quality for the new date (2022-05-18)
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44699
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 100112043
$ws.Range("G3").Value = "Pepino dulce"
$ws.Range("H3").Value = "Cultivar IV Región"
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 13000
$ws.Range("N3").Value = "$/bandeja 18 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 722
$ws.Range("Q3").Value = 18
$ws.Range("R3").Value = "Hortaliza"
